# MasterExecutor_Sanity.xlsx - "TC19 Updated for ECTEST"
# Update RunMode (column E) values to proper case (Yes/No) and flip most
# rows to "Yes" (i.e. enabled) while TC19 and its neighboring rows stay "No"
# (disabled) and become hidden via the autofilter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that remain disabled ("No") - TC19_Verify_ShippingMethod and the
# surrounding checkout-related test cases (rows 12, 19-22).
$noRows = @(12, 19, 20, 21, 22)

for ($r = 2; $r -le 32; $r++) {
    if ($noRows -contains $r) {
        $ws.Range("E$r").Value2 = "No"
    } else {
        $ws.Range("E$r").Value2 = "Yes"
    }
}

# Hide the disabled rows (as the autofilter would when RunMode != Yes).
foreach ($r in $noRows) {
    $ws.Rows($r).Hidden = $true
}

# Update the view: scroll a bit further down and select the header row.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("A1:XFD1").Select()

$wb.Save()
